$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.659.13"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.78%  "
$ws.Range("D3").Value = "'3.183.26"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.68%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'589.19"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.36%  "
$ws.Range("D6").Value = "'135.97"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.52%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'3.180.69"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.76%  "
$ws.Range("E9").Value = "  -2.68%  "
$ws.Range("E10").Value = "  -4.76%  "
$ws.Range("D11").Value = "'5.28"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.02%  "
$ws.Range("E12").Value = "  -3.82%  "
$ws.Range("E13").Value = "  -4.80%  "
$ws.Range("D14").Value = "'33.42"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.74%  "
$ws.Range("D15").Value = "'3.704.93"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.85%  "
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").Value = "'3.177.25"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.97%  "
$ws.Range("D18").Value = "'62.633.69"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.98%  "
$ws.Range("D19").Value = "'6.55"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.91%  "
$ws.Range("D20").Value = "'456.79"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.07%  "
$ws.Range("D21").Value = "'13.93"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").Value = "'0.704"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.89%  "
$ws.Range("D24").Value = "'13.41"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "'83.48"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -2.75%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'6.89"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -6.32%  "
$ws.Range("E30").Value = "  -4.52%  "
$ws.Range("E31").Value = "  -6.53%  "
$ws.Range("D32").Value = "'27.41"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.50%  "
$ws.Range("D33").Value = "'0.104"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("E34").Value = "  -5.65%  "
$ws.Range("D35").Value = "'1.04"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.78%  "
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("D37").Value = "'51.05"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.53%  "
$ws.Range("D38").Value = "'0.0₃0700"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -6.48%  "
$ws.Range("E39").Value = "  -4.08%  "
$ws.Range("D40").Value = "'409.59"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.56%  "
$ws.Range("D41").Value = "'2.72"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("D42").Value = "'2.853.77"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.30%  "
$ws.Range("D43").Value = "'8.01"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.62%  "
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").Value = "'36.50"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.15%  "
$ws.Range("E46").Value = "  -6.09%  "
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D49").Value = "'125.10"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").Value = "'25.52"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.01%  "
$ws.Range("E51").Value = "  -3.65%  "
